$wb = $excel.ActiveWorkbook

# Select the full sheet on Slovakia (mirrors the "select all cells" state
# left behind on the source sheet after duplicating it for the new market).
$src = $wb.Worksheets.Item("Slovakia")
$src.Cells.Select()

# Duplicate the Slovakia sheet and place the copy right after it, at the
# end of the tab strip.
$src.Copy($null, $src)

# The duplicated sheet becomes the active sheet, positioned right after
# "Slovakia".
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "Italy"

# Update the market-specific values for Italy.
$newSheet.Range("B2").Value = "Italy Market"
$newSheet.Range("B4").Value = "NGC-3145/T2454/T2453"

# Leave the cursor on B4, matching the state captured after typing the
# new part number.
$newSheet.Range("B4").Select()
